$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (UnitMass) values for the "+ loading" table (rows 2-21)
# and the "- loading" table (rows 23-42), per the commit diff.
$ws.Cells.Item(2, 3).Value = 23
$ws.Cells.Item(3, 3).Value = 39
$ws.Cells.Item(4, 3).Value = 167
$ws.Cells.Item(5, 3).Value = 79
$ws.Cells.Item(6, 3).Value = 137
$ws.Cells.Item(7, 3).Value = 46
$ws.Cells.Item(8, 3).Value = 111
$ws.Cells.Item(9, 3).Value = 98
$ws.Cells.Item(10, 3).Value = 45
$ws.Cells.Item(11, 3).Value = 67
$ws.Cells.Item(13, 3).Value = 66
$ws.Cells.Item(14, 3).Value = 105
$ws.Cells.Item(15, 3).Value = 327
$ws.Cells.Item(16, 3).Value = 119
$ws.Cells.Item(17, 3).Value = 55
$ws.Cells.Item(18, 3).Value = 199
$ws.Cells.Item(19, 3).Value = 122
$ws.Cells.Item(20, 3).Value = 71
$ws.Cells.Item(21, 3).Value = 43
$ws.Cells.Item(23, 3).Value = 27
$ws.Cells.Item(24, 3).Value = 56
$ws.Cells.Item(25, 3).Value = 40
$ws.Cells.Item(26, 3).Value = 44
$ws.Cells.Item(27, 3).Value = 60
$ws.Cells.Item(28, 3).Value = 87
$ws.Cells.Item(29, 3).Value = 24
$ws.Cells.Item(30, 3).Value = 59
$ws.Cells.Item(31, 3).Value = 61
$ws.Cells.Item(32, 3).Value = 73
$ws.Cells.Item(33, 3).Value = 100
$ws.Cells.Item(34, 3).Value = 58
$ws.Cells.Item(35, 3).Value = 72
$ws.Cells.Item(36, 3).Value = 32
$ws.Cells.Item(37, 3).Value = 20
$ws.Cells.Item(38, 3).Value = 101
$ws.Cells.Item(39, 3).Value = 70
$ws.Cells.Item(40, 3).Value = 75
$ws.Cells.Item(41, 3).Value = 28
$ws.Cells.Item(42, 3).Value = 18
